# Apply the "merger / adding_ltp in main_scanner" update:
#  - relabel the two header groups from "45 days"/"15 days" to "45_days"/"15_days"
#  - refresh the ltp (last traded price) column (B) with newly scanned values
#  - move the active selection to P12 (reflects where the user left off editing)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header labels (row 1 merged headers) ---
$ws.Range("C1").Value = "45_days"
$ws.Range("L1").Value = "15_days"

# --- Refreshed ltp values (column B) ---
$ws.Range("B3").Value = 601.85
$ws.Range("B4").Value = 4890
$ws.Range("B5").Value = 24340.4
$ws.Range("B6").Value = 173.25
$ws.Range("B7").Value = 225.4

# --- Active selection moves to P12 ---
$ws.Range("P12").Select()
